# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.911.28"
$ws.Range("E2").Value = "  -4.43%  "
$ws.Range("D3").Value = "2.461.43"
$ws.Range("E3").Value = "  -5.79%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'548.07"
$ws.Range("E5").Value = "  -4.49%  "
$ws.Range("D6").Value = "'145.57"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -3.75%  "
$ws.Range("D9").Value = "2.460.57"
$ws.Range("E9").Value = "  -5.70%  "
$ws.Range("E10").Value = "  -9.14%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = "  -7.95%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -7.82%  "
$ws.Range("D14").Value = "'26.03"
$ws.Range("E14").Value = "  -7.56%  "
$ws.Range("D15").Value = "2.907.98"
$ws.Range("E15").Value = "  -5.59%  "
$ws.Range("D16").Value = "'0.0000163"
$ws.Range("E16").Value = "  -9.39%  "
$ws.Range("D17").Value = "60.829.35"
$ws.Range("E17").Value = "  -4.38%  "
$ws.Range("D18").Value = "2.462.09"
$ws.Range("E18").Value = "  -5.40%  "
$ws.Range("D19").Value = "'11.04"
$ws.Range("E19").Value = "  -8.00%  "
$ws.Range("D20").Value = "'6.92"
$ws.Range("E20").Value = "  -8.16%  "
$ws.Range("D21").Value = "'4.17"
$ws.Range("E21").Value = "  -7.70%  "
$ws.Range("D22").Value = "'318.84"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'63.37"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "0.0₃0978"
$ws.Range("E26").Value = "  -9.31%  "
$ws.Range("D27").Value = "2.582.35"
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'533.96"
$ws.Range("E29").Value = "  -9.26%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.47"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").Value = "'8.29"
$ws.Range("E31").Value = "  -9.09%  "
$ws.Range("D32").Value = "'7.60"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -6.56%  "
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("E35").Value = "  -9.31%  "
$ws.Range("D36").Value = "'5.86"
$ws.Range("E36").Value = "  -11.61%  "
$ws.Range("D37").Value = "'4.85"
$ws.Range("E37").Value = "  -10.37%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.375"
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("D40").Value = "'18.35"
$ws.Range("E40").Value = "  -6.85%  "
$ws.Range("D41").Value = "'145.87"
$ws.Range("E41").Value = "  -5.35%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'1.70"
$ws.Range("E43").Value = "  -9.30%  "
$ws.Range("D44").Value = "'39.88"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").Value = "'2.28"
$ws.Range("E45").Value = "  -8.88%  "
$ws.Range("D46").Value = "'146.93"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("D47").Value = "'3.56"
$ws.Range("E47").Value = "  -8.31%  "
$ws.Range("D48").Value = "'20.72"
$ws.Range("E48").Value = "  -12.80%  "
$ws.Range("D49").Value = "'0.0530"
$ws.Range("E49").Value = "  -10.03%  "
$ws.Range("D50").Value = "'0.584"
$ws.Range("E50").Value = "  -7.31%  "
$ws.Range("D51").Value = "'0.0938"
$ws.Range("E51").Value = "  -6.17%  "
